# Update agv_logs_M2 sheet: recompute cumulative timestamps (column A)
# and reassign job ids (column D) for rows 2-18, then remove the last
# data row (row 19) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values: row number -> (timestamp, job_id)
$updates = @{
    2  = @(154.1,              "J12")
    3  = @(704.7,               "J21")
    4  = @(876.9000000000001,   "J23")
    5  = @(1401.5,               "J9")
    6  = @(1653.899999999999,   "J16")
    7  = @(1719.999999999999,   "J22")
    8  = @(2112.499999999999,   "J19")
    9  = @(2284.699999999999,    "J8")
    10 = @(3075.699999999998,   "J24")
    11 = @(3201.999999999998,   "J24")
    12 = @(3282.199999999998,   "J14")
    13 = @(3987.299999999997,   "J30")
    14 = @(4388.099999999997,   "J11")
    15 = @(4548.599999999999,   "J11")
    16 = @(4563.599999999999,   "J14")
    17 = @(4627,                "J15")
    18 = @(4709.1,               "J6")
}

foreach ($r in $updates.Keys) {
    $pair = $updates[$r]
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("D$r").Value = $pair[1]
}

# Remove the 19th data row entirely so the sheet's dimension shrinks to F18
$ws.Rows.Item(19).Delete()
